# "changed BIGINT to INT in Code Slide"
#
# Slide 11 ("SQL-Code, Höchstzahlverfahren") holds the SQL source code
# inside the first shape ("Rechteck 3"). One of the lines declares the
# "seats" column as BIGINT; the author narrowed it to INT.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$hit = $tr.Find("BIGINT")
if ($hit -ne $null) {
    $hit.Text = "INT"
}

# The slide also carried a leftover, completely empty content
# placeholder ("Inhaltsplatzhalter 2") that the author cleared out
# while touching this slide.
$ph = $s.Shapes.Item("Inhaltsplatzhalter 2")
if ($ph -ne $null) {
    $ph.Delete()
}
